## Insert a new, separate run containing a single space " " immediately
## before the existing "Plan: " run at the very start of the document
## (same run formatting: Times New Roman, 12pt, en-GB).
##
## A plain InsertBefore()/TypeText() would just splice the space into the
## neighbouring run's <w:t> (since the formatting is identical), but the
## target edit keeps it as its own <w:r> element, so we build the run via
## a tiny WordprocessingML fragment and insert it with Range.InsertXML,
## which preserves run boundaries exactly as authored.

$d = $word.ActiveDocument

# Locate the "Plan:" run robustly (it's the very first text in the body).
$target = $d.Content.Duplicate
$target.Find.Execute("Plan:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Collapse to the insertion point right before "Plan:".
$target.Collapse(1)

$runXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
<w:lang w:val="en-GB"/>
</w:rPr>
<w:t xml:space="preserve"> </w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($runXml)
